{"js": "// Rewrite the intro paragraph: \"paper\" -> \"article\", \"Life Expectancy\" ->\n// \"life expectancy\", restructure the two numbered clauses, and swap\n// \"Delimits\" for \"extends\" (plus related wording tweaks), per the commit\n// \"updates to app text\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\n\nconst newText =\n  \"Our article highlights key disparities in life expectancy in US MSAs \" +\n  \"with > 1 million people. This app does two key things: 1) translates \" +\n  \"the article figures into interactive visualizations, allowing more \" +\n  \"granular exploration of the article analysis, and 2) extends \" +\n  \"visualizations to all US MSAs including those with populations of \" +\n  \"less than 1 million residents.\";\n\n// Replace the entire paragraph's text in one shot so the original run's\n// (lack of) formatting is preserved rather than introducing new runs.\nparagraph.insertText(newText, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Rewrite the intro paragraph: \"paper\" -> \"article\", \"Life Expectancy\" ->\n# \"life expectancy\", restructure the two numbered clauses, and swap\n# \"Delimits\" for \"extends\" (plus related wording tweaks), per the commit\n# \"updates to app text\".\n$d = $word.ActiveDocument\n\n$newText = \"Our article highlights key disparities in life expectancy in US MSAs with > 1 million people. This app does two key things: 1) translates the article figures into interactive visualizations, allowing more granular exploration of the article analysis, and 2) extends visualizations to all US MSAs including those with populations of less than 1 million residents.\"\n\n# Replace the entire first paragraph's text in one shot so the original\n# run's (lack of) formatting is preserved rather than introducing new runs.\n$p = $d.Paragraphs(1)\n$p.Range.Text = $newText\n"}
